# modificando ejemplos TEI 0.2.1
# Fill in the (previously empty) "Definition" column (D) for each concept
# row on the "Concepts" sheet with the same text as the "Display" column
# (C), for rows 2 through 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

for ($r = 2; $r -le 8; $r++) {
    $display = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 4).Value = $display
}
